# Reorder the data rows (rows 2-22) on both worksheets according to the
# drag-and-drop reorder performed by the user. The set of rows/values is
# unchanged - only their vertical position moved.
#
# Mapping: for each NEW row position (2..22), this says which OLD row
# position (2..22) supplied its data.
$perm = @(5, 9, 15, 7, 20, 12, 18, 8, 13, 17, 2, 21, 3, 6, 4, 19, 16, 22, 14, 10, 11)

$wb = $excel.ActiveWorkbook

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    # Column C ("normalized_upc") holds digit-strings that must stay text
    # (not be auto-coerced to numbers) when we rewrite them below.
    $ws.Range("C2:C22").NumberFormat = "@"

    # Snapshot all data rows (A2:H22) before we start overwriting anything.
    $snapshot = $ws.Range("A2:H22").Value2

    # Build the reordered 2-D array (1-based, rows 1..21, cols 1..8).
    $reordered = New-Object 'object[,]' 21,8
    for ($i = 1; $i -le 21; $i++) {
        $srcRow = $perm[$i - 1] - 1   # old sheet row -> snapshot row index (1-based)
        for ($j = 1; $j -le 8; $j++) {
            $reordered[$i - 1, $j - 1] = $snapshot[$srcRow, $j]
        }
    }

    $ws.Range("A2:H22").Value2 = $reordered
}
